$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 143, shifting existing rows 143:165 down to 144:166
$ws.Rows("143:143").Insert()

# Populate the newly inserted row 143 with the new record
$ws.Cells.Item(143, 1).Value = 11
$ws.Cells.Item(143, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(143, 3).Value = "Bíobío"
$ws.Cells.Item(143, 4).Value = 44504
$ws.Cells.Item(143, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(143, 5).Value = 8
$ws.Cells.Item(143, 6).Value = 100112008
$ws.Cells.Item(143, 7).Value = "Coliflor"
$ws.Cells.Item(143, 8).Value = "Sin especificar"
$ws.Cells.Item(143, 9).Value = "Primera"
$ws.Cells.Item(143, 10).Value = 2700
$ws.Cells.Item(143, 11).Value = 650
$ws.Cells.Item(143, 12).Value = 700
$ws.Cells.Item(143, 13).Value = 672
$ws.Cells.Item(143, 14).Value = "$/unidad"
$ws.Cells.Item(143, 15).Value = "Región Metropolitana"
$ws.Cells.Item(143, 16).Value = 672
$ws.Cells.Item(143, 17).Value = 1
$ws.Cells.Item(143, 18).Value = "Hortaliza"
